$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$dCell = $ws.Range("D2")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '26.649.06'
$dCell.Style = $dStyle
$ws.Range("E2").Value = '  +0.63%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$dCell = $ws.Range("D3")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '1.847.50'
$dCell.Style = $dStyle
$ws.Range("E3").Value = '  -0.01%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$dCell = $ws.Range("D4")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.9997'
$dCell.Style = $dStyle
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$dCell = $ws.Range("D5")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '262.88'
$dCell.Style = $dStyle
$ws.Range("E5").Value = '  -0.94%  '

$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$dCell = $ws.Range("D6")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.9998'
$dCell.Style = $dStyle
$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$dCell = $ws.Range("D7")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.5330'
$dCell.Style = $dStyle
$ws.Range("E7").Value = '  +2.57%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$dCell = $ws.Range("D8")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.3216'
$dCell.Style = $dStyle
$ws.Range("E8").Value = '  -2.02%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$dCell = $ws.Range("D9")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.06936'
$dCell.Style = $dStyle
$ws.Range("E9").Value = '  +1.82%  '

$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$dCell = $ws.Range("D10")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '19.18'
$dCell.Style = $dStyle
$ws.Range("E10").Value = '  +1.36%  '

$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$dCell = $ws.Range("D11")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.7843'
$dCell.Style = $dStyle
$ws.Range("E11").Value = '  +0.52%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$dCell = $ws.Range("D12")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.07823'
$dCell.Style = $dStyle
$ws.Range("E12").Value = '  +0.77%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$dCell = $ws.Range("D13")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '1.829.32'
$dCell.Style = $dStyle
$ws.Range("E13").Value = '  -1.03%  '

$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$dCell = $ws.Range("D14")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '89.28'
$dCell.Style = $dStyle
$ws.Range("E14").Value = '  +1.28%  '

$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$dCell = $ws.Range("D15")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '5.053'
$dCell.Style = $dStyle
$ws.Range("E15").Value = '  +0.63%  '

$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$dCell = $ws.Range("D16")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '14.18'
$dCell.Style = $dStyle
$ws.Range("E16").Value = '  +1.59%  '

$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$dCell = $ws.Range("D17")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.9994'
$dCell.Style = $dStyle
$ws.Range("E17").Value = '  +0.14%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$dCell = $ws.Range("D18")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.000007997'
$dCell.Style = $dStyle
$ws.Range("E18").Value = '  +0.05%  '

$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$dCell = $ws.Range("D19")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.9998'
$dCell.Style = $dStyle
$ws.Range("E19").Value = '  +0.02%  '

$ws.Range("B20").Value = 'WrappedBTC'
$ws.Range("C20").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$dCell = $ws.Range("D20")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '26.671.58'
$dCell.Style = $dStyle
$ws.Range("E20").Value = '  +0.62%  '

$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$dCell = $ws.Range("D21")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '2.073.34'
$dCell.Style = $dStyle
$ws.Range("E21").Value = '  +0.02%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$dCell = $ws.Range("D22")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '4.652'
$dCell.Style = $dStyle
$ws.Range("E22").Value = '  +0.34%  '

$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$dCell = $ws.Range("D23")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '6.033'
$dCell.Style = $dStyle
$ws.Range("E23").Value = '  +0.54%  '

$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$dCell = $ws.Range("D24")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '9.418'
$dCell.Style = $dStyle
$ws.Range("E24").Value = '  -1.88%  '

$ws.Range("B25").Value = 'LidoDAOToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$dCell = $ws.Range("D25")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '2.231'
$dCell.Style = $dStyle
$ws.Range("E25").Value = '  +1.69%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$dCell = $ws.Range("D26")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '142.73'
$dCell.Style = $dStyle
$ws.Range("E26").Value = '  -1.03%  '

$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$dCell = $ws.Range("D27")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '1.700'
$dCell.Style = $dStyle
$ws.Range("E27").Value = '  +1.94%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$dCell = $ws.Range("D28")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '17.13'
$dCell.Style = $dStyle
$ws.Range("E28").Value = '  +0.46%  '

$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$dCell = $ws.Range("D29")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '111.70'
$dCell.Style = $dStyle
$ws.Range("E29").Value = '  -0.47%  '

$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$dCell = $ws.Range("D30")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '4.280'
$dCell.Style = $dStyle
$ws.Range("E30").Value = '  +2.39%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$dCell = $ws.Range("D31")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.08757'
$dCell.Style = $dStyle
$ws.Range("E31").Value = '  +0.33%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$dCell = $ws.Range("D32")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '4.119'
$dCell.Style = $dStyle
$ws.Range("E32").Value = '  -0.47%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$dCell = $ws.Range("D33")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.04881'
$dCell.Style = $dStyle
$ws.Range("E33").Value = '  +0.80%  '

$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$dCell = $ws.Range("D34")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.7388'
$dCell.Style = $dStyle
$ws.Range("E34").Value = '  +1.82%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$dCell = $ws.Range("D35")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '1.144'
$dCell.Style = $dStyle
$ws.Range("E35").Value = '  +0.78%  '

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$dCell = $ws.Range("D36")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '2.864'
$dCell.Style = $dStyle
$ws.Range("E36").Value = '  +0.45%  '

$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$dCell = $ws.Range("D37")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '3.111'
$dCell.Style = $dStyle
$ws.Range("E37").Value = '  +0.18%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$dCell = $ws.Range("D38")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '2.354'
$dCell.Style = $dStyle
$ws.Range("E38").Value = '  +5.81%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$dCell = $ws.Range("D39")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.01745'
$dCell.Style = $dStyle
$ws.Range("E39").Value = '  -2.40%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$dCell = $ws.Range("D40")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.4848'
$dCell.Style = $dStyle
$ws.Range("E40").Value = '  -1.03%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$dCell = $ws.Range("D41")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.9087'
$dCell.Style = $dStyle
$ws.Range("E41").Value = '  -0.81%  '

$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$dCell = $ws.Range("D42")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '109.75'
$dCell.Style = $dStyle
$ws.Range("E42").Value = '  -1.53%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$dCell = $ws.Range("D43")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '5.915'
$dCell.Style = $dStyle
$ws.Range("E43").Value = '  -2.93%  '

$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$dCell = $ws.Range("D44")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.9999'
$dCell.Style = $dStyle
$ws.Range("E44").Value = '  +0.08%  '

$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$dCell = $ws.Range("D45")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '7.749'
$dCell.Style = $dStyle
$ws.Range("E45").Value = '  -0.26%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$dCell = $ws.Range("D46")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.4231'
$dCell.Style = $dStyle
$ws.Range("E46").Value = '  +0.80%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$dCell = $ws.Range("D47")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.1258'
$dCell.Style = $dStyle
$ws.Range("E47").Value = '  +0.74%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$dCell = $ws.Range("D48")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '9.090'
$dCell.Style = $dStyle
$ws.Range("E48").Value = '  -0.04%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$dCell = $ws.Range("D49")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '35.09'
$dCell.Style = $dStyle
$ws.Range("E49").Value = '  +0.11%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$dCell = $ws.Range("D50")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.05841'
$dCell.Style = $dStyle
$ws.Range("E50").Value = '  -1.72%  '

$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$dCell = $ws.Range("D51")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = '0.8974'
$dCell.Style = $dStyle
$ws.Range("E51").Value = '  +0.97%  '
